$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# --- Update the "Desc Body" text in B5 (shorten it; move the benchmark
#     detail sentence out into a new row) -----------------------------------
$ws.Range("B5").Value = "The number of homeless persons was over 116,000 in 2016 – a 29.8 per cent increase from just under 90,000 in 2006 and a 13.7 per cent increase from 102,000 in 2011. The national benchmark was not met."
$ws.Rows.Item(5).RowHeight = 41.75

# Row 6 ("Most of the increase...") keeps its text but gets a new height.
$ws.Rows.Item(6).RowHeight = 41.75

# --- Insert a new row 7 holding the new explanatory paragraph -------------
$ws.Rows.Item(7).Insert()
$ws.Range("B7").Value = "State and Territory assessments are made against the agreed performance indicator — proportion of Australians who are homeless. This was around 50 homeless for every 10,000 people nationally in 2016"
$ws.Rows.Item(7).RowHeight = 63

# Keep the selection on B5 like in the edited workbook.
$ws.Range("B5").Select()
